# Update the "Förändrad" (Changed) date column (C) for rows 2-28
# from 45502 (2024-07-29) to 45503 (2024-07-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45502) {
        $cell.Value2 = 45503
    }
}
